# Update Pooh Points site
#  - Status column (G) text for in-progress games changed from the old
#    clock-based label to "End of 2nd Half"
#  - Status column (G) is also narrowed by one character

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# Rows in the Players sheet whose status (column G) reads "13:48 - 2nd Half"
$rows = @(3, 10, 11, 16, 18, 25, 28, 31, 34, 38, 41, 42, 43, 48, 52, 53, 54, 55, 57, 58, 61)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "13:48 - 2nd Half") {
        $cell.Value = "End of 2nd Half"
    }
}

# Column G width: 18 -> 17 (ColumnWidth is offset from the stored sheet
# width by the default cell-padding constant, ~0.8333333333333333)
$ws.Columns.Item(7).ColumnWidth = 17 - 0.8333333333333333
